# Update the "取得日時" (retrieved datetime) column for rows 2-14 on the
# "ランサーズ" sheet from 2026-02-11 18:59:45 to 2026-02-11 19:08:48.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-02-11 19:08:48"

for ($row = 2; $row -le 14; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
